# Renamed repo / fixed output folder path: the underlying MAG file listing
# changed, which removed three MAG entries from this sheet
# (even_MAG-GUT2323.fa, even_MAG-GUT27781.fa, even_MAG-GUT78207.fa) and
# shifted the remaining rows up to fill the gaps. The sheet's used range
# shrinks from A1:G24 down to A1:G21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (even_MAG-GUT2323.fa) and the original row 10 (even_MAG-GUT27781.fa)
# are dropped. Deleting row 9 twice removes both, because after the first
# delete the old row 10 shifts up into row 9.
$ws.Range("A9").EntireRow.Delete()
$ws.Range("A9").EntireRow.Delete()

# The original last data row, row 24 (even_MAG-GUT78207.fa), is also
# dropped. After removing the two rows above, it now sits at row 22.
$ws.Range("A22").EntireRow.Delete()
